$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(6)

$shp.TextFrame.TextRange.Text = "WHY IT MATTERS:"

$shp.Left = -229.15212598425197
$shp.Top = 450.8511811023622
$shp.Width = 1094.3896062992126
$shp.Height = 60.75
